# Forest data update - 2026-02-06 12:30
#
# Two listings that were sitting in the "New" sheet (Gulbene / Ludza) are
# promoted to the bottom of "Previously added" (rows 471-472), and five
# freshly scraped listings (Liepaja, Limbazi, Madona, Rezekne, Valka) take
# over the "New" sheet.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# =======================================================================
# PART 1 - append the two previously-"New" rows to "Previously added"
# =======================================================================

# Copy the formatting of the last existing row (470) down onto the two
# new rows (471-472) first, so every cell already carries the right
# style before we touch values.
$wsPrev.Range("A470:F470").Copy()
$wsPrev.Range("A471:F472").PasteSpecial(-4122)

# --- row 471 (ex "New"!row 2 - Gulbene listing) ---
$wsPrev.Range("A471").Value = "https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/gulbene/ggxmf.html"
$wsPrev.Range("B471").Value = "22 000 €"
$wsPrev.Range("C471").Value = "Gulbene un raj."
$wsPrev.Range("D471").Value = "2 ha."
$wsPrev.Range("E471").NumberFormat = "@"
$wsPrev.Range("E471").Value = "50440140001"
$wsPrev.Range("F471").Value = 46058.45763888889

$wsPrev.Hyperlinks.Add($wsPrev.Range("A471"), "https://www.ss.com/msg/lv/real-estate/wood/gulbene-and-reg/gulbene/ggxmf.html")

# --- row 472 (ex "New"!row 3 - Ludza listing) ---
$wsPrev.Range("A472").Value = "https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/zalesjes-pag/bmhjim.html"
$wsPrev.Range("B472").Value = "8 000 €"
$wsPrev.Range("C472").Value = "Ludza un raj."
$wsPrev.Range("D472").Value = "11 ha."
$wsPrev.Range("E472").NumberFormat = "@"
$wsPrev.Range("E472").Value = "68960050098"
$wsPrev.Range("F472").Value = 46058.49097222222

$wsPrev.Hyperlinks.Add($wsPrev.Range("A472"), "https://www.ss.com/msg/lv/real-estate/wood/ludza-and-reg/zalesjes-pag/bmhjim.html")

# Re-apply the template formatting once more so that adding the
# hyperlinks (which stamps its own "Hyperlink" style) and forcing a text
# number format on column E don't leave the wrong style behind.
$wsPrev.Range("A470:F470").Copy()
$wsPrev.Range("A471:F472").PasteSpecial(-4122)

# =======================================================================
# PART 2 - replace the "New" sheet's 2 rows with 5 freshly scraped rows
# =======================================================================

# Make rows 2-6 all carry the same formatting as the current row 2/3
# template before writing any values.
$wsNew.Range("A2:F2").Copy()
$wsNew.Range("A2:F6").PasteSpecial(-4122)

# --- row 2 - Liepaja listing ---
$wsNew.Range("A2").Value = "https://www.ss.com/msg/lv/real-estate/wood/liepaja-and-reg/dunikas-pag/bljip.html"
$wsNew.Range("B2").Value = "6 400 €"
$wsNew.Range("C2").Value = "Liepāja un raj."
$wsNew.Range("D2").Value = "2 ha."
$wsNew.Range("E2").NumberFormat = "@"
$wsNew.Range("E2").Value = "64520050033"
$wsNew.Range("F2").Value = 46059.478472222225

# --- row 3 - Limbazi listing ---
$wsNew.Range("A3").Value = "https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/katvaru-pag/blfig.html"
$wsNew.Range("B3").Value = "45 000 €"
$wsNew.Range("C3").Value = "Limbaži un raj."
$wsNew.Range("D3").Value = "5 ha."
$wsNew.Range("E3").Value = "6652 006 0004"
$wsNew.Range("F3").Value = 46059.490277777775

# --- row 4 - Madona listing ---
$wsNew.Range("A4").Value = "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/sarkanu-pag/idpic.html"
$wsNew.Range("B4").Value = "45 000 €"
$wsNew.Range("C4").Value = "Madona un raj."
$wsNew.Range("D4").Value = "5 ha."
$wsNew.Range("E4").NumberFormat = "@"
$wsNew.Range("E4").Value = "70900080055"
$wsNew.Range("F4").Value = 46058.90486111111

# --- row 5 - Rezekne listing ---
$wsNew.Range("A5").Value = "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/sakstagala-pag/fdlph.html"
$wsNew.Range("B5").Value = "10 000 €"
$wsNew.Range("C5").Value = "Rēzekne un raj."
$wsNew.Range("D5").Value = "3 ha."
$wsNew.Range("E5").NumberFormat = "@"
$wsNew.Range("E5").Value = "78860020157"
$wsNew.Range("F5").Value = 46058.62291666667

# --- row 6 - Valka listing ---
$wsNew.Range("A6").Value = "https://www.ss.com/msg/lv/real-estate/wood/valka-and-reg/planu-pag/ooflx.html"
$wsNew.Range("B6").Value = "100 000 €"
$wsNew.Range("C6").Value = "Valka un raj."
$wsNew.Range("D6").Value = "4 ha."
$wsNew.Range("E6").NumberFormat = "@"
$wsNew.Range("E6").Value = "94760120100"
$wsNew.Range("F6").Value = 46059.51666666666

# Hyperlinks for the 5 new rows of "New" (added after the text values so
# the displayed text is the scraped link itself, matching row 1's links).
$wsNew.Hyperlinks.Add($wsNew.Range("A2"), "https://www.ss.com/msg/lv/real-estate/wood/liepaja-and-reg/dunikas-pag/bljip.html")
$wsNew.Hyperlinks.Add($wsNew.Range("A3"), "https://www.ss.com/msg/lv/real-estate/wood/limbadzi-and-reg/katvaru-pag/blfig.html")
$wsNew.Hyperlinks.Add($wsNew.Range("A4"), "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/sarkanu-pag/idpic.html")
$wsNew.Hyperlinks.Add($wsNew.Range("A5"), "https://www.ss.com/msg/lv/real-estate/wood/rezekne-and-reg/sakstagala-pag/fdlph.html")
$wsNew.Hyperlinks.Add($wsNew.Range("A6"), "https://www.ss.com/msg/lv/real-estate/wood/valka-and-reg/planu-pag/ooflx.html")

# Re-apply formatting one final time so every cell ends up with the
# original column style (hyperlink-add/text-format tweaks above can
# otherwise leave a stray style behind).
$wsNew.Range("A2:F2").Copy()
$wsNew.Range("A2:F6").PasteSpecial(-4122)

Write-Host "Forest data updated: +2 rows in 'Previously added', 'New' sheet refreshed with 5 listings."
